$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 607.8570999999999
$ws.Range("I12").Value = 607.8570999999999
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 607.8570999999999
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -437.8570999999999
$ws.Range("H18").Value = 1404.5
$ws.Range("I18").Value = 1404.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1404.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1120.5
$ws.Range("H40").Value = 6350.7
$ws.Range("I40").Value = 5358.2856
$ws.Range("J40").Value = 8666.333000000001
$ws.Range("K40").Value = 5358.2856
$ws.Range("L40").Value = 8666.333000000001
$ws.Range("M40").Value = -5183.2856
$ws.Range("N40").Value = -9016.333000000001
$ws.Range("H68").Value = 75281.5
$ws.Range("I68").Value = 100268
$ws.Range("J68").Value = 50295
$ws.Range("K68").Value = 100268
$ws.Range("L68").Value = 50295
$ws.Range("M68").Value = -99519
$ws.Range("N68").Value = -51793
$ws.Range("H71").Value = 75281.5
$ws.Range("I71").Value = 100268
$ws.Range("J71").Value = 50295
$ws.Range("K71").Value = 300804
$ws.Range("L71").Value = 150885
$ws.Range("M71").Value = -297060
$ws.Range("N71").Value = -158373
$ws.Range("H74").Value = 3114.5
$ws.Range("I74").Value = 3114.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3114.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2178.5
$ws.Range("H77").Value = 3114.5
$ws.Range("I77").Value = 3114.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 15572.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -10892.5
$ws.Range("H96").Value = 746.2857
$ws.Range("I96").Value = 105
$ws.Range("J96").Value = 2349.5
$ws.Range("K96").Value = 315
$ws.Range("L96").Value = 7048.5
$ws.Range("M96").Value = 1058
$ws.Range("N96").Value = -9794.5
$ws.Range("H98").Value = 479.66666
$ws.Range("I98").Value = 275.6
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 275.6
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 1222.4
$ws.Range("N98").Value = -4496
$ws.Range("H113").Value = 5196
$ws.Range("I113").Value = 5247.125
$ws.Range("J113").Value = 4991.5
$ws.Range("K113").Value = 5247.125
$ws.Range("L113").Value = 4991.5
$ws.Range("M113").Value = -1993.125
$ws.Range("N113").Value = -11499.5
$ws.Range("H122").Value = 479.66666
$ws.Range("I122").Value = 275.6
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 826.8000000000001
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = 1623.2
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 4584.2666
$ws.Range("I132").Value = 1147
$ws.Range("J132").Value = 18333.334
$ws.Range("K132").Value = 3441
$ws.Range("L132").Value = 55000.00199999999
$ws.Range("M132").Value = -911
$ws.Range("N132").Value = -60060.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 42.5
$ws.Range("I5").Value = 23.333334
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 23.333334
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 88.66666599999999
$ws.Range("N5").Value = -324
$ws.Range("H96").Value = 10022222
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 10022222
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 10022222
$ws.Range("N96").Value = -10027714
$ws.Range("H102").Value = 4195.0586
$ws.Range("I102").Value = 2027.8182
$ws.Range("J102").Value = 8168.3335
$ws.Range("K102").Value = 2027.8182
$ws.Range("L102").Value = 8168.3335
$ws.Range("M102").Value = -405.8181999999999
$ws.Range("N102").Value = -11412.3335
$ws.Range("H122").Value = 2359.8
$ws.Range("I122").Value = 2487.25
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 7461.75
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -5011.75
$ws.Range("N122").Value = -10450
$ws.Range("H132").Value = 18742
$ws.Range("I132").Value = 18742
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 56226
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -53696

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 42.5
$ws.Range("I4").Value = 23.333334
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 23.333334
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 91.66666599999999
$ws.Range("N4").Value = -330
$ws.Range("H42").Value = 220000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 220000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 220000
$ws.Range("N42").Value = -220656
$ws.Range("H86").Value = 4809.4707
$ws.Range("I86").Value = 2997
$ws.Range("J86").Value = 6420.5557
$ws.Range("K86").Value = 2997
$ws.Range("L86").Value = 6420.5557
$ws.Range("M86").Value = -1874
$ws.Range("N86").Value = -8666.555700000001
$ws.Range("H89").Value = 4809.4707
$ws.Range("I89").Value = 2997
$ws.Range("J89").Value = 6420.5557
$ws.Range("K89").Value = 14985
$ws.Range("L89").Value = 32102.7785
$ws.Range("M89").Value = -9369
$ws.Range("N89").Value = -43334.7785
$ws.Range("H99").Value = 2995
$ws.Range("I99").Value = 2325.6667
$ws.Range("J99").Value = 3999
$ws.Range("K99").Value = 2325.6667
$ws.Range("L99").Value = 3999
$ws.Range("M99").Value = -827.6667000000002
$ws.Range("N99").Value = -6995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1686
$ws.Range("I22").Value = 1471.2
$ws.Range("J22").Value = 2223
$ws.Range("K22").Value = 1471.2
$ws.Range("L22").Value = 2223
$ws.Range("M22").Value = -1121.2
$ws.Range("N22").Value = -2923
$ws.Range("H94").Value = 4645.4287
$ws.Range("I94").Value = 626
$ws.Range("J94").Value = 10004.667
$ws.Range("K94").Value = 626
$ws.Range("L94").Value = 10004.667
$ws.Range("M94").Value = -175
$ws.Range("N94").Value = -10906.667
$ws.Range("H103").Value = 17192.4
$ws.Range("I103").Value = 17192.4
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 17192.4
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -16020.4
$ws.Range("H107").Value = 1019.5
$ws.Range("I107").Value = 409.8
$ws.Range("J107").Value = 1629.2
$ws.Range("K107").Value = 409.8
$ws.Range("L107").Value = 1629.2
$ws.Range("M107").Value = 1510.2
$ws.Range("N107").Value = -5469.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 330.33334
$ws.Range("I44").Value = 168.25
$ws.Range("J44").Value = 460
$ws.Range("K44").Value = 504.75
$ws.Range("L44").Value = 1380
$ws.Range("M44").Value = -106.75
$ws.Range("N44").Value = -2176
$ws.Range("H107").Value = 276.92856
$ws.Range("I107").Value = 222.6
$ws.Range("J107").Value = 412.75
$ws.Range("K107").Value = 667.8
$ws.Range("L107").Value = 1238.25
$ws.Range("M107").Value = 1252.2
$ws.Range("N107").Value = -5078.25
$ws.Range("H128").Value = 810000
$ws.Range("I128").Value = 810000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 2430000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -2425020
$ws.Range("H138").Value = 8210
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 8210
$ws.Range("K138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("M138").Value = 24630
$ws.Range("N138").Value = -34910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3129.3
$ws.Range("I102").Value = 2397.5715
$ws.Range("J102").Value = 4836.6665
$ws.Range("K102").Value = 2397.5715
$ws.Range("L102").Value = 4836.6665
$ws.Range("M102").Value = -775.5715
$ws.Range("N102").Value = -8080.6665
$ws.Range("H126").Value = 5582.2
$ws.Range("I126").Value = 4303.6665
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 12910.9995
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -10440.9995
$ws.Range("N126").Value = -27440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 11749.5
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 11749.5
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 11749.5
$ws.Range("N2").Value = -11973.5
$ws.Range("H16").Value = 2500.5
$ws.Range("I16").Value = 2500.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2500.5
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2330.5
$ws.Range("H46").Value = 5998.125
$ws.Range("I46").Value = 4400
$ws.Range("J46").Value = 6724.5454
$ws.Range("K46").Value = 4400
$ws.Range("L46").Value = 6724.5454
$ws.Range("M46").Value = -4212
$ws.Range("N46").Value = -7100.5454
$ws.Range("H122").Value = 4899.6
$ws.Range("I122").Value = 4874.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 14623.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -12173.5
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 53000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 53000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 53000
$ws.Range("N64").Value = -53496
$ws.Range("H67").Value = 53000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 53000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 53000
$ws.Range("N67").Value = -54716
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H102").Value = 46999.75
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 46999.75
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 46999.75
$ws.Range("N102").Value = -53489.75
$ws.Range("H107").Value = 615.0833
$ws.Range("I107").Value = 615.0833
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1845.2499
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 74.75009999999997
$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 7500
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -12400
